# Correction in SA (simulated annealing) algorithm and run_19 log.
# Updates the Fitness column (C) for generations 0-72 (rows 2-74)
# to reflect the corrected simulated-annealing fitness trace.
# Rows 75+ (generation 73 onward) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value  = 8423
$ws.Range("C11:C19").Value = 7945
$ws.Range("C20:C33").Value = 7812
$ws.Range("C34:C44").Value = 7754
$ws.Range("C45:C56").Value = 7721
$ws.Range("C57:C73").Value = 7704
$ws.Range("C74").Value     = 7343
